$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column labels replacing old Test_Case / Run_Mode
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "runmode"

# Data row (row 2): sample credentials/run-mode values
$ws.Range("A2").Value = "pramodnp.pnp@gmail.com"
$ws.Range("B2").Value = "pramodnp1995"
$ws.Range("C2").Value = "y"

# Column widths for A and B
$ws.Columns.Item(1).ColumnWidth = 31.33203125
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

# Highlight header row with a yellow fill
$ws.Range("A1:C1").Interior.Color = 65535

# Move the active selection to D8
$ws.Range("D8").Select()
